# Weekly crime-data refresh: shift report week from 1/1/2024-1/7/2024
# to 1/8/2024-1/14/2024 (Volume 31 Number 1 -> Number 2) and update the
# precinct crime-complaint table (rows 14-30) with the new weeks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header rich-text runs ---
$ws.Range("A8").Characters(21, 1).Text = "2"
$ws.Range("C9").Characters(27, 8).Text = "1/8/2024"
$ws.Range("C9").Characters(46, 8).Text = "1/14/2024"

# --- Cells that change from a blank-marker shared string ("0" / "***.*")
#     to a real number: clone formatting+type from a stable donor cell
#     that keeps the same marker text before AND after this edit, then
#     overwrite with the new number. ---
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))

# --- Cells that change from a real number to the blank-marker shared
#     string: same clone-from-donor trick, using the opposite donor. ---

# --- Plain numeric updates (values only; existing cell style already matches) ---
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -85.714285714285
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 66.666666666666
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -73.684210526315
$ws.Range("N16").Value = -88.636363636363
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -42.857142857142
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 29.166666666666
$ws.Range("I17").Value = 15
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 36.363636363636
$ws.Range("L17").Value = 15.384615384615
$ws.Range("M17").Value = 36.363636363636
$ws.Range("N17").Value = -48.275862068965
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 166.666666666667
$ws.Range("K18").Value = 100
$ws.Range("L18").Value = -77.777777777777
$ws.Range("M18").Value = -85.714285714285
$ws.Range("N18").Value = -95.833333333333
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -56.756756756756
$ws.Range("I19").Value = 7
$ws.Range("J19").Value = 15
$ws.Range("K19").Value = -53.333333333333
$ws.Range("L19").Value = -56.25
$ws.Range("M19").Value = -56.25
$ws.Range("N19").Value = -91.358024691358
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = 71.428571428571
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -83.098591549295
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = -5.208333333333
$ws.Range("I21").Value = 42
$ws.Range("J21").Value = 38
$ws.Range("K21").Value = 10.526315789473
$ws.Range("L21").Value = -16
$ws.Range("M21").Value = -45.454545454545
$ws.Range("N21").Value = -85.159010600706
$ws.Range("M23").Value = -50
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 2.247191011235
$ws.Range("I24").Value = 46
$ws.Range("J24").Value = 39
$ws.Range("K24").Value = 17.948717948717
$ws.Range("L24").Value = -4.166666666666
$ws.Range("M24").Value = 12.195121951219
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -35.714285714285
$ws.Range("F25").Value = 44
$ws.Range("H25").Value = 12.820512820512
$ws.Range("I25").Value = 17
$ws.Range("J25").Value = 24
$ws.Range("K25").Value = -29.166666666666
$ws.Range("L25").Value = -10.526315789473
$ws.Range("M25").Value = -32
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("L26").Value = -50
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 3
$ws.Range("L27").Value = 200
$ws.Range("F28").Value = 5
$ws.Range("F29").Value = 2

# --- Cells that change from the blank-marker shared string to a real number:
#     set the value then restore the columns usual number format so the
#     cell lands back on the normal "count" / "% change" style. ---
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M14").Value = -100
$ws.Range("M14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I15").Value = 1
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = "#,##0"
$ws.Range("K15").Value = 0
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("I18").Value = 2
$ws.Range("I18").NumberFormat = "#,##0"
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("F23").Value = 1
$ws.Range("F23").NumberFormat = "#,##0"
$ws.Range("I23").Value = 1
$ws.Range("I23").NumberFormat = "#,##0"
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I26").Value = 1
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 1
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("K26").Value = 0
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D27").Value = 4
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -75
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J27").Value = 4
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("K27").Value = -25
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M28").Value = -100
$ws.Range("M28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M29").Value = -100
$ws.Range("M29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L30").Value = -100
$ws.Range("L30").NumberFormat = '#,##0.0;"-"#,##0.0'
